$wb = $excel.ActiveWorkbook

# --- "Edit Repayment Schedule" sheet: update the selected cell ---
$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsEdit.Activate() | Out-Null
$wsEdit.Range("A6").Select() | Out-Null

# --- "Repayment schedule" sheet: insert a new (blank) column before N ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Activate() | Out-Null
$wsRepay.Columns.Item(14).Insert() | Out-Null
# the inserted column keeps the width of the column immediately to its left
$wsRepay.Columns.Item(14).ColumnWidth = $wsRepay.Columns.Item(13).ColumnWidth
$wsRepay.Range("S10").Select() | Out-Null
